# Update cryptos list with new price / volume(1h) figures.
# Values that look numeric are written with a leading apostrophe (forces
# text entry, matching the source data's text-typed Price/Volume columns)
# and the cell style is reset to "Normal" afterwards so no stray
# quote-prefix formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows where both Price (D) and Volume(1h) (E) change ---
$ws.Range("D2").Value = "41.542.79"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "2.194.49"
$ws.Range("E3").Value = "  -2.61%  "

$ws.Range("D5").Value = "'228.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").Value = "'0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.57%  "

$ws.Range("D7").Value = "'59.75"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.04%  "

$ws.Range("D9").Value = "'0.399"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.97%  "

$ws.Range("D10").Value = "'56.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.95%  "

$ws.Range("D11").Value = "'0.0881"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.67%  "

$ws.Range("D13").Value = "2.526.80"
$ws.Range("E13").Value = "  -2.34%  "

$ws.Range("D14").Value = "'15.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.76%  "

$ws.Range("D15").Value = "'22.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.67%  "

$ws.Range("D16").Value = "'5.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.87%  "

$ws.Range("D17").Value = "'0.788"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.20%  "

$ws.Range("D18").Value = "2.215.83"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").Value = "41.499.32"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D22").Value = "'6.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.10%  "

$ws.Range("D23").Value = "'241.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.36%  "

$ws.Range("D24").Value = "'0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D27").Value = "'9.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.20%  "

$ws.Range("D28").Value = "'168.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("D30").Value = "'1.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.22%  "

$ws.Range("D31").Value = "'19.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.18%  "

$ws.Range("D32").Value = "'2.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.96%  "

$ws.Range("D34").Value = "'4.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.01%  "

$ws.Range("D35").Value = "'4.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.90%  "

$ws.Range("D36").Value = "'0.0643"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("D37").Value = "'6.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.73%  "

$ws.Range("D38").Value = "'2.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.08%  "

$ws.Range("D39").Value = "'3.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.17%  "

$ws.Range("D40").Value = "'0.000238"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.63%  "

$ws.Range("D42").Value = "'0.0235"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.85%  "

$ws.Range("D43").Value = "'8.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.78%  "

$ws.Range("D44").Value = "'0.0950"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.86%  "

$ws.Range("D47").Value = "1.458.88"
$ws.Range("E47").Value = "  -3.53%  "

$ws.Range("D48").Value = "'4.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -17.95%  "

$ws.Range("D49").Value = "'16.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.02%  "

$ws.Range("D50").Value = "'2.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("D51").Value = "'1.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.09%  "

# --- Rows where only Volume(1h) (E) changes ---
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("E29").Value = "  -7.11%  "
$ws.Range("E33").Value = "  -4.03%  "
$ws.Range("E41").Value = "  +0.04%  "

# --- Rows 20 / 21 swap: Litecoin moves up to rank 20, ShibaInu moves down to rank 21 ---
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'71.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.01%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0" + [char]0x2083 + "0894"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.65%  "

# --- Rows 45 / 46 swap: Aave moves up to rank 45, TrustWalletToken moves down to rank 46 ---
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'96.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.20%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.47%  "
